$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the " ; " separator with " - " in the cells that use it,
# as described in the commit message (";" conflicts with the cell separator).
$ws.Range("B10").Value = "CDTSFILIA - "
$ws.Range("B13").Value = "CCODDFT - `nCXLAPTZ - `nCYLAPTZ"
$ws.Range("B18").Value = "CDTSCVT - `nCDTMCVT"
$ws.Range("A29").Value = "ILTASIT - `nILTASEU"
$ws.Range("C34").Value = "TA_SEUIL - `nTA_INFOS_SEUIL"
$ws.Range("D34").Value = "DATE_SAISIE - `nDATE_SAISIE"

# Update the view state: scroll back to top-left and move selection to E19.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E19").Select()
